# Applies the "REF_GenerateYearlyReport -> REF_GenerateYearlyReport_Performer"
# rebrand edits to the Config workbook (Settings sheet gains new rows, asset
# path value changes, column B widened, selection moved).

$wb = $excel.ActiveWorkbook
$settings = $wb.Worksheets.Item("Settings")

# B5: process / asset name changes from REF_GenerateYearlyReport to the
# "-Performer" variant.
$settings.Range("B5").Value = "REF_GenerateYearlyReport-Performer"

# Rows 11-14: four brand-new Name/Value/Description entries describing the
# work item status and the new report URLs used by the Performer process.
# Cells are written in the same order the original author entered them so
# that the workbook's shared-string table is rebuilt in the same sequence.
$settings.Range("A11").Value = "Status"
$settings.Range("B11").Value = "Completed"
$settings.Range("B12").Value = "https://acme-test.uipath.com/work-items"
$settings.Range("B13").Value = "https://acme-test.uipath.com/reports/download"
$settings.Range("A14").Value = "UploadYearlyReport_URL"
$settings.Range("A12").Value = "WorkItems_URL"
$settings.Range("A13").Value = "DownloadMonthlyReport_URL"
$settings.Range("B14").Value = "https://acme-test.uipath.com/reports/upload"
$settings.Range("C11").Value = "Comment Status"
$settings.Range("C13").Value = "Download Monthly Report URL"
$settings.Range("C14").Value = "Upload Yearly Report URL"
$settings.Range("C12").Value = "Work Items url"

# B9: ReportDirPath asset value now points at the new Performer project's
# Reports folder instead of the old Output folder.
$settings.Range("B9").Value = "C:\Users\Charlie\OneDrive\Documents\UiPath Advance\REF_GenerateYearlyReport_Performer\Data\Reports\"

# Column B on the Settings sheet is widened to fit the longer new URL/path
# text that was just added (target stored width is 101 characters; the
# runtime's ColumnWidth setter adds a constant ~0.8333 padding on save, so
# back that offset out to land exactly on 101 in the saved file).
$settings.Columns.Item(2).ColumnWidth = 100.16666666666667

# The author's last selection before saving moved to B8.
$settings.Range("B8").Select()

$wb.Save()
